# Normalize the "Recorded By" (column G) values: when the comma-separated
# list of recorder names includes an exact "System" entry, move it so that
# the first and last entries of the list are swapped (this puts "System"
# at the front in the common two-entry case, and otherwise keeps the rest
# of the relative ordering while rotating the extremes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range to know how many rows contain data.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value()

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -ge 2 -and ($parts -contains "System")) {
        $first = $parts[0]
        $last = $parts[$parts.Length - 1]
        $parts[0] = $last
        $parts[$parts.Length - 1] = $first
        $newText = [string]::Join(", ", $parts)
        $cell.Value = $newText
    }
}
